# Rename the two logo pictures that live in the document's headers/footers.
#
#   - The BTec logo (in both headers) goes from "image1.jpg" -> "image2.jpg"
#   - The Pearson logo (in both footers) goes from "image2.png" -> "image1.png"
#
# Both wdHeaderFooterPrimary (1) and wdHeaderFooterFirstPage (2) slots carry
# their own copy of the picture, so every header/footer combination is
# visited and updated.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-LogoInStory($story, $newName) {
    if (-not $story.Exists) {
        return
    }
    $rng = $story.Range
    if ($rng.InlineShapes.Count -lt 1) {
        return
    }
    # Re-fetch the inline shape through its own Range before writing - doing
    # the rename straight off the header/footer Range's collection can hand
    # back a handle that the host considers stale.
    $probe = $rng.InlineShapes.Item(1)
    $shp = $probe.Range.InlineShapes.Item(1)
    $shp.Name = $newName
}

# Headers: BTec_Logo-Orange, image1.jpg -> image2.jpg
Rename-LogoInStory $sec.Headers.Item(1) "image2.jpg"
Rename-LogoInStory $sec.Headers.Item(2) "image2.jpg"

# Footers: PearsonLogo.png, image2.png -> image1.png
Rename-LogoInStory $sec.Footers.Item(1) "image1.png"
Rename-LogoInStory $sec.Footers.Item(2) "image1.png"
